# "amend errors: data entered in the wrong columns, labels around the
# wrong way"
#
# Row 11 (Schvartz (1977) / thermo) had n_exp / mu_exp / sd_exp entered
# one column to the left of where they belong:
#   C11 (n_exp)  held 53.1  (that's mu_exp) -> should be 21
#   D11 (mu_exp) held 1.15  (that's sd_exp) -> should be 53.1
#   E11 (sd_exp) held 21    (that's n_exp)  -> should be 1.15
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C11").Value = 21
$ws.Range("D11").Value = 53.1
$ws.Range("E11").Value = 1.15

# The author's selection when they saved moved from D7 to G9.
$ws.Range("G9").Select()
